$wb = $excel.ActiveWorkbook

# This script applies a scheduled-runner style refresh of computed
# market/profit columns (H:N) across multiple sheets, matching the
# authoritative diff. Values are pure numeric overwrites; a couple of
# rows also gain/lose a cell (M122 added on ARM/CRP, N16 removed on LTW).

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 23902.139
$ws.Range("I28").Value = 39159.25
$ws.Range("J28").Value = 5124.154
$ws.Range("K28").Value = 39159.25
$ws.Range("L28").Value = 5124.154
$ws.Range("M28").Value = -38674.25
$ws.Range("N28").Value = -6094.154
$ws.Range("H86").Value = 3289.276
$ws.Range("I86").Value = 1521
$ws.Range("J86").Value = 4219.9473
$ws.Range("K86").Value = 1521
$ws.Range("L86").Value = 4219.9473
$ws.Range("M86").Value = -398
$ws.Range("N86").Value = -6465.9473
$ws.Range("H89").Value = 3289.276
$ws.Range("I89").Value = 1521
$ws.Range("J89").Value = 4219.9473
$ws.Range("K89").Value = 7605
$ws.Range("L89").Value = 21099.7365
$ws.Range("M89").Value = -1989
$ws.Range("N89").Value = -32331.7365
$ws.Range("H113").Value = 2282.8572
$ws.Range("I113").Value = 1639.4615
$ws.Range("J113").Value = 2840.4666
$ws.Range("K113").Value = 1639.4615
$ws.Range("L113").Value = 2840.4666
$ws.Range("M113").Value = 1614.5385
$ws.Range("N113").Value = -9348.4666
$ws.Range("H132").Value = 3311.5425
$ws.Range("I132").Value = 3357.4092
$ws.Range("J132").Value = 3203.4285
$ws.Range("K132").Value = 10072.2276
$ws.Range("L132").Value = 9610.2855
$ws.Range("M132").Value = -7542.2276
$ws.Range("N132").Value = -14670.2855
$ws.Range("H137").Value = 3447.2642
$ws.Range("I137").Value = 1219.8
$ws.Range("J137").Value = 4797.242
$ws.Range("K137").Value = 3659.4
$ws.Range("L137").Value = 14391.726
$ws.Range("M137").Value = -1109.4
$ws.Range("N137").Value = -19491.726
$ws.Range("H138").Value = 2456.9124
$ws.Range("I138").Value = 1209.8889
$ws.Range("J138").Value = 3579.2334
$ws.Range("K138").Value = 3629.6667
$ws.Range("L138").Value = 10737.7002
$ws.Range("M138").Value = 1510.3333
$ws.Range("N138").Value = -21017.7002
$ws.Range("H141").Value = 4390.6113
$ws.Range("I141").Value = 1441.9333
$ws.Range("J141").Value = 19134
$ws.Range("K141").Value = 4325.7999
$ws.Range("L141").Value = 57402
$ws.Range("M141").Value = 854.2001
$ws.Range("N141").Value = -67762

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5860.557
$ws.Range("I32").Value = 4873.5166
$ws.Range("K32").Value = 4873.5166
$ws.Range("M32").Value = -4586.5166
$ws.Range("H45").Value = 1115.9678
$ws.Range("I45").Value = 956.3043
$ws.Range("J45").Value = 1575
$ws.Range("K45").Value = 956.3043
$ws.Range("L45").Value = 1575
$ws.Range("M45").Value = -579.3043
$ws.Range("N45").Value = -2329
$ws.Range("H110").Value = 1043.2069
$ws.Range("I110").Value = 737.4666999999999
$ws.Range("J110").Value = 1370.7858
$ws.Range("K110").Value = 737.4666999999999
$ws.Range("L110").Value = 1370.7858
$ws.Range("M110").Value = 1307.5333
$ws.Range("N110").Value = -5460.7858
$ws.Range("H122").Value = 1040.2106
$ws.Range("I122").Value = 896.0909
$ws.Range("J122").Value = 1238.375
$ws.Range("K122").Value = 2688.2727
$ws.Range("L122").Value = 3715.125
$ws.Range("N122").Value = -8615.125
$ws.Range("H132").Value = 7908.9355
$ws.Range("I132").Value = 5555.625
$ws.Range("K132").Value = 16666.875
$ws.Range("M132").Value = -14136.875
$ws.Range("M122").Value = -238.2727

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27687.352
$ws.Range("I20").Value = 591.6667
$ws.Range("J20").Value = 63250.438
$ws.Range("K20").Value = 591.6667
$ws.Range("L20").Value = 63250.438
$ws.Range("M20").Value = -344.6667
$ws.Range("N20").Value = -63744.438
$ws.Range("H22").Value = 135.8
$ws.Range("I22").Value = 94.75
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 94.75
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 78.25
$ws.Range("N22").Value = -646
$ws.Range("H68").Value = 46500.5
$ws.Range("J68").Value = 46500.5
$ws.Range("L68").Value = 46500.5
$ws.Range("N68").Value = -48122.5
$ws.Range("H71").Value = 46500.5
$ws.Range("J71").Value = 46500.5
$ws.Range("L71").Value = 139501.5
$ws.Range("N71").Value = -147613.5
$ws.Range("H107").Value = 1474.2174
$ws.Range("I107").Value = 784.7692
$ws.Range("J107").Value = 2370.5
$ws.Range("K107").Value = 784.7692
$ws.Range("L107").Value = 2370.5
$ws.Range("M107").Value = 1135.2308
$ws.Range("N107").Value = -6210.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 903.64703
$ws.Range("I58").Value = 713.5897
$ws.Range("J58").Value = 1521.3334
$ws.Range("K58").Value = 713.5897
$ws.Range("L58").Value = 1521.3334
$ws.Range("M58").Value = -510.5897
$ws.Range("N58").Value = -1927.3334
$ws.Range("H62").Value = 4160.972
$ws.Range("J62").Value = 2830.9375
$ws.Range("L62").Value = 2830.9375
$ws.Range("N62").Value = -4078.9375
$ws.Range("H65").Value = 4160.972
$ws.Range("J65").Value = 2830.9375
$ws.Range("L65").Value = 14154.6875
$ws.Range("N65").Value = -20394.6875
$ws.Range("H122").Value = 1539.25
$ws.Range("I122").Value = 1385.6666
$ws.Range("K122").Value = 4156.9998
$ws.Range("H132").Value = 12197652
$ws.Range("I132").Value = 16396375
$ws.Range("J132").Value = 1359.238
$ws.Range("K132").Value = 49189125
$ws.Range("L132").Value = 4077.714
$ws.Range("M132").Value = -49186595
$ws.Range("N132").Value = -9137.714
$ws.Range("H134").Value = 2487.145
$ws.Range("I134").Value = 2701.6833
$ws.Range("J134").Value = 1056.8889
$ws.Range("K134").Value = 8105.0499
$ws.Range("L134").Value = 3170.6667
$ws.Range("M134").Value = -5570.0499
$ws.Range("N134").Value = -8240.6667
$ws.Range("H136").Value = 903.64703
$ws.Range("I136").Value = 713.5897
$ws.Range("J136").Value = 1521.3334
$ws.Range("K136").Value = 2140.7691
$ws.Range("L136").Value = 4564.0002
$ws.Range("M136").Value = 409.2309
$ws.Range("N136").Value = -9664.0002
$ws.Range("M122").Value = -1706.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2728.5715
$ws.Range("J55").Value = 3600
$ws.Range("L55").Value = 10800
$ws.Range("N55").Value = -11154
$ws.Range("H132").Value = 826.8
$ws.Range("I132").Value = 826.8
$ws.Range("K132").Value = 7441.2
$ws.Range("M132").Value = -4911.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 42647.5
$ws.Range("J68").Value = 42647.5
$ws.Range("L68").Value = 42647.5
$ws.Range("N68").Value = -44269.5
$ws.Range("H69").Value = 48000
$ws.Range("J69").Value = 48000
$ws.Range("L69").Value = 48000
$ws.Range("N69").Value = -49498
$ws.Range("H71").Value = 42647.5
$ws.Range("J71").Value = 42647.5
$ws.Range("L71").Value = 127942.5
$ws.Range("N71").Value = -136054.5
$ws.Range("H72").Value = 48000
$ws.Range("J72").Value = 48000
$ws.Range("L72").Value = 144000
$ws.Range("N72").Value = -151488
$ws.Range("H80").Value = 3287.9
$ws.Range("I80").Value = 2205
$ws.Range("J80").Value = 3752
$ws.Range("K80").Value = 2205
$ws.Range("L80").Value = 3752
$ws.Range("M80").Value = -1207
$ws.Range("N80").Value = -5748
$ws.Range("H83").Value = 3287.9
$ws.Range("I83").Value = 2205
$ws.Range("J83").Value = 3752
$ws.Range("K83").Value = 11025
$ws.Range("L83").Value = 18760
$ws.Range("M83").Value = -6033
$ws.Range("N83").Value = -28744
$ws.Range("H122").Value = 9706.666999999999
$ws.Range("I122").Value = 10316.363
$ws.Range("K122").Value = 30949.089
$ws.Range("M122").Value = -28499.089
$ws.Range("H126").Value = 1606.909
$ws.Range("I126").Value = 1374.4
$ws.Range("J126").Value = 1800.6666
$ws.Range("K126").Value = 4123.200000000001
$ws.Range("L126").Value = 5401.9998
$ws.Range("M126").Value = -1653.200000000001
$ws.Range("N126").Value = -10341.9998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 330.52942
$ws.Range("I16").Value = 330.52942
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 330.52942
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -160.52942
$ws.Range("H40").Value = 1599.9231
$ws.Range("I40").Value = 1350.4
$ws.Range("J40").Value = 2431.6667
$ws.Range("K40").Value = 1350.4
$ws.Range("L40").Value = 2431.6667
$ws.Range("M40").Value = -1214.4
$ws.Range("N40").Value = -2703.6667
$ws.Range("H122").Value = 55594.316
$ws.Range("I122").Value = 79531.69500000001
$ws.Range("J122").Value = 3730
$ws.Range("K122").Value = 238595.085
$ws.Range("L122").Value = 11190
$ws.Range("M122").Value = -236145.085
$ws.Range("N122").Value = -16090
$ws.Range("H127").Value = 15000
$ws.Range("J127").Value = 15000
$ws.Range("L127").Value = 15000
$ws.Range("N127").Value = -24920

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 242.3
$ws.Range("I107").Value = 236.23529
$ws.Range("J107").Value = 276.66666
$ws.Range("K107").Value = 708.70587
$ws.Range("L107").Value = 829.9999799999999
$ws.Range("M107").Value = 1211.29413
$ws.Range("N107").Value = -4669.99998
$ws.Range("H122").Value = 12502221
$ws.Range("I122").Value = 22223844
$ws.Range("J122").Value = 2991.4285
$ws.Range("K122").Value = 66671532
$ws.Range("L122").Value = 8974.2855
$ws.Range("M122").Value = -66669082
$ws.Range("N122").Value = -13874.2855
$ws.Range("H126").Value = 768
$ws.Range("I126").Value = 519.17645
$ws.Range("J126").Value = 1191
$ws.Range("K126").Value = 1557.52935
$ws.Range("L126").Value = 3573
$ws.Range("M126").Value = 912.4706499999998
$ws.Range("N126").Value = -8513
$ws.Range("H132").Value = 2051.5898
$ws.Range("I132").Value = 2373.5964
$ws.Range("J132").Value = 1177.5714
$ws.Range("K132").Value = 7120.789199999999
$ws.Range("L132").Value = 3532.7142
$ws.Range("M132").Value = -4590.789199999999
$ws.Range("N132").Value = -8592.7142
$ws.Range("H136").Value = 2680.9036
$ws.Range("I136").Value = 3315.7144
$ws.Range("J136").Value = 1766.0294
$ws.Range("K136").Value = 9947.143199999999
$ws.Range("L136").Value = 5298.0882
$ws.Range("M136").Value = -7397.143199999999
$ws.Range("N136").Value = -10398.0882

# LTW row 16: N16 column value is removed entirely (leve has no HQ profit figure)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N16").ClearContents()
